$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.714.92"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "'2.618.08"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'594.88"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").Value = "'149.39"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.108"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.381"
$ws.Range("E10").Value = "  +3.52%  "
$ws.Range("D11").Value = "'5.58"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D13").Value = "'27.43"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "'3.089.49"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "'63.585.60"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "'2.613.01"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "'12.11"
$ws.Range("E18").Value = "  +5.61%  "
$ws.Range("D19").Value = "'4.60"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").Value = "'347.75"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").Value = "'6.84"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").Value = "'66.09"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "'1.72"
$ws.Range("E25").Value = "  +11.90%  "
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "'9.16"
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "'8.06"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").Value = "'543.82"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'2.02"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "'0.0₃0845"
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("D34").Value = "'1.74"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "'5.22"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").Value = "'168.19"
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.405"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("D40").Value = "'19.30"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'168.60"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "'39.82"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "'3.90"
$ws.Range("E44").Value = "  +3.41%  "
$ws.Range("D45").Value = "'0.0586"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").Value = "'21.32"
$ws.Range("E46").Value = "  -5.22%  "
$ws.Range("D47").Value = "'0.626"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("E48").Value = "  +11.26%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").Value = "'0.0963"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "'19.04"
$ws.Range("E51").Value = "  +1.31%  "
